$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.709.21'
$ws.Range("E2").Value = '  +2.47%  '

$ws.Range("D3").Value = '3.391.14'
$ws.Range("E3").Value = '  +1.42%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.79'
$ws.Range("E5").Value = '  +1.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.16'
$ws.Range("E6").Value = '  +3.32%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '3.390.37'
$ws.Range("E8").Value = '  +1.35%  '

$ws.Range("E9").Value = '  +0.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.46'
$ws.Range("E10").Value = '  +0.01%  '

$ws.Range("E11").Value = '  +5.85%  '

$ws.Range("E12").Value = '  +3.82%  '

$ws.Range("D13").Value = '3.971.98'
$ws.Range("E13").Value = '  +1.60%  '

$ws.Range("E14").Value = '  +2.37%  '

$ws.Range("D16").Value = '3.395.12'
$ws.Range("E16").Value = '  +2.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.35'
$ws.Range("E17").Value = '  +1.88%  '

$ws.Range("D18").Value = '61.737.56'
$ws.Range("E18").Value = '  +2.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.20'
$ws.Range("E19").Value = '  +4.40%  '

$ws.Range("E20").Value = '  +2.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.45'
$ws.Range("E21").Value = '  +1.98%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '387.13'
$ws.Range("E22").Value = '  +8.82%  '

$ws.Range("E23").Value = '  +1.29%  '

$ws.Range("D24").Value = '3.537.83'
$ws.Range("E24").Value = '  +1.86%  '

$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("E26").Value = '  +13.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.18'
$ws.Range("E27").Value = '  +2.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.65'
$ws.Range("E28").Value = '  +2.05%  '

$ws.Range("E29").Value = '  -3.82%  '

$ws.Range("E30").Value = '  +0.23%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.25'
$ws.Range("E31").Value = '  +3.32%  '

$ws.Range("E32").Value = '  +4.08%  '

$ws.Range("E33").Value = '  +1.52%  '

$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("D35").Value = '3.423.81'
$ws.Range("E35").Value = '  +1.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.47'
$ws.Range("E36").Value = '  +2.23%  '

$ws.Range("E37").Value = '  +0.44%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.96'
$ws.Range("E38").Value = '  +0.89%  '

$ws.Range("E39").Value = '  +2.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '163.22'
$ws.Range("E40").Value = '  +2.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0786'
$ws.Range("E41").Value = '  +1.59%  '

$ws.Range("E42").Value = '  +12.14%  '

$ws.Range("E43").Value = '  +4.34%  '

$ws.Range("E44").Value = '  +0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.22'
$ws.Range("E45").Value = '  +2.72%  '

$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.83'
$ws.Range("E46").Value = '  +2.21%  '

$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.44'
$ws.Range("E47").Value = '  +0.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.62'
$ws.Range("E48").Value = '  +4.18%  '

$ws.Range("E49").Value = '  +1.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.24'
$ws.Range("E50").Value = '  +3.01%  '

$ws.Range("D51").Value = '2.354.47'
$ws.Range("E51").Value = '  +7.98%  '
